# Correção da prova 2.
#
# The grader/error explanation paragraph wrongly blamed an unidentified
# syntax error "on line 41" ("... na linha 41 porém nem eu nem o monitor
# conseguimos identificar."). The correction pins down the real cause: the
# error is actually on line 35, caused by using an 'else' with an argument
# where an 'else if' was required.

$d = $word.ActiveDocument

# --- Main content fix -------------------------------------------------
# Replace the unresolved "linha 41 ... não conseguimos identificar" remark
# with the accurate diagnosis pointing at line 35 and its real cause.
$d.Content.Find.Execute(
    "Na linha 41 utilizo um ‘=’ no lugar de um sinal de menos ‘-’. O compilador informa um erro de sintaxe na linha 41 porém nem eu nem o monitor conseguimos identificar.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Na linha 41 utilizo um ‘=’ no lugar de um sinal de menos ‘-’. O compilador informa um erro de sintaxe na linha 35 ocasionado pela utilização de um ‘else’ com argumento, deveria ser utilizado um ‘else if’.",
    2) | Out-Null

# --- Cosmetic run clean-up (no visible text change) --------------------
# These two Find/Replace no-ops just coalesce runs that were previously
# split mid-word ("...P" + "2", "ve" + "c" + "tor"), matching how Word
# re-flows runs after any edit touches a paragraph.
$d.Content.Find.Execute(
    "Estrutura de Dados 1 – Correção da questão 13 – P2",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Estrutura de Dados 1 – Correção da questão 13 – P2",
    2) | Out-Null

$d.Content.Find.Execute(
    "primeiro verifico se  um dos vector já está vazio, se sim,",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "primeiro verifico se  um dos vector já está vazio, se sim,",
    2) | Out-Null

Write-Output "edit applied"
